$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.174419283866882
$ws.Range("B1").Value = 2.250952482223511
$ws.Range("C1").Value = 4.702385902404785
$ws.Range("D1").Value = 2.631445169448853
$ws.Range("E1").Value = 1.226363182067871
